$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the "Good Morning" greeting to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the selection change recorded in the saved view state
$ws.Range("E8").Select()
